$d = $word.ActiveDocument

# The document has a full set of tracked insertions/deletions (all by
# "Rodney Thayer") covering wording clarifications throughout the test
# description, plus a couple of paragraph-mark insert/delete markers.
# The target revision simply accepts every one of these tracked changes
# (keeping the inserted text, discarding the deleted text, and merging
# any paragraphs whose paragraph mark was itself deleted) -- i.e. a
# plain "Accept All Changes in Document".
$d.AcceptAllRevisions()
